$d = $word.ActiveDocument
$d.Content.Find.Execute("JSON/XML scripting, in depth testing and debugging", $true, $false, $false, $false, $false, $true, 1, $false, "JSON/XML scripting, data analysis in R, testing and debugging", 2)
